# UndoRedoSequenceDiagram.pptx - "Update UG DG diagrams to reflect Event Manager"
#
# Renames the AddressBook-era class names in the sequence diagram on slide 1
# to their Event Manager equivalents, and nudges a few shapes that had to be
# resized/repositioned to fit the new (longer/shorter) labels.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "Rectangle 65" (id=56) - big rounded rectangle behind the right-hand
#    lifeline group. Only repositioned (same size).
# ---------------------------------------------------------------------------
$rect65 = Get-ShapeById $s 56
$rect65.Left = 507.544094488189
$rect65.Top  = 28.98763779527559

# ---------------------------------------------------------------------------
# 2) "Rectangle 62" (id=16) - ":Address" / "BookParser" lifeline header
#    becomes ":EventManager" / "Parser", and shrinks from 16pt to 12pt.
# ---------------------------------------------------------------------------
$rect62a = Get-ShapeById $s 16
$rect62a.TextFrame.TextRange.Text = ":EventManager" + [char]13 + "Parser"
$rect62a.TextFrame.TextRange.Font.Size = 12

# ---------------------------------------------------------------------------
# 3) "TextBox 78" (id=79) - "undoAddressBook()" becomes "undoEventManager()"
#    (only the middle run's text changes, colors/err flags stay put), text
#    shrinks from 12pt to 10.5pt, and the box is shortened back to a single
#    line's height.
# ---------------------------------------------------------------------------
$textbox78 = Get-ShapeById $s 79
$full78 = $textbox78.TextFrame.TextRange
$mid78 = $full78.Characters(5, 11)   # "AddressBook"
$mid78.Text = "EventManager"
$textbox78.TextFrame.TextRange.Font.Size = 10.5
$textbox78.Height = 12.723070866141732

# ---------------------------------------------------------------------------
# 4) "Rectangle 62" (id=84) - ":VersionedAddressBook" becomes
#    ":VersionedEventManager", shrinking from 16pt to 14pt.
# ---------------------------------------------------------------------------
$rect62b = Get-ShapeById $s 84
$full84 = $rect62b.TextFrame.TextRange
$run84 = $full84.Characters(2, 20)   # "VersionedAddressBook"
$run84.Text = "VersionedEventManager"
$rect62b.TextFrame.TextRange.Font.Size = 14

# ---------------------------------------------------------------------------
# 5) "TextBox 87" (id=88) - "resetData(ReadOnlyAddressBook)" becomes
#    "resetData(ReadOnlyEventManager)"; box widens to fit the longer label.
# ---------------------------------------------------------------------------
$textbox87 = Get-ShapeById $s 88
$full87 = $textbox87.TextFrame.TextRange
$mid87 = $full87.Characters(11, 19)  # "ReadOnlyAddressBook"
$mid87.Text = "ReadOnlyEventManager"
$textbox87.Left   = 690.514094488189
$textbox87.Top    = 257.27064592125987
$textbox87.Width  = 180.67968503937007
$textbox87.Height = 14.540629921259843
